$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: turn the A-column-only shared-formula demo into a two column
# Value / Formula table (A = input numbers, B = A*n shared formula).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A1").Value = "Value"
$ws1.Range("B1").Value = "Formula"

for ($r = 2; $r -le 19; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}

$ws1.Range("B2").Formula = "=A2*10"
$ws1.Range("B3:B19").Formula = "=A3*10"

# Margins matching the legacy Excel 2003/2007 defaults used by the target file.
$ws1.PageSetup.LeftMargin = 0.75 * 72
$ws1.PageSetup.RightMargin = 0.75 * 72
$ws1.PageSetup.TopMargin = 1 * 72
$ws1.PageSetup.BottomMargin = 1 * 72
$ws1.PageSetup.HeaderMargin = 0.5 * 72
$ws1.PageSetup.FooterMargin = 0.5 * 72

# Selection moves to B4 (the first "dragged-down" shared formula cell).
$ws1.Range("B4").Select()

# ---------------------------------------------------------------------
# Sheet2 / Sheet3: same margin normalization, nothing else changes.
# ---------------------------------------------------------------------
foreach ($name in @("Sheet2", "Sheet3")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.PageSetup.LeftMargin = 0.75 * 72
    $ws.PageSetup.RightMargin = 0.75 * 72
    $ws.PageSetup.TopMargin = 1 * 72
    $ws.PageSetup.BottomMargin = 1 * 72
    $ws.PageSetup.HeaderMargin = 0.5 * 72
    $ws.PageSetup.FooterMargin = 0.5 * 72
}

# ---------------------------------------------------------------------
# Workbook default font: Calibri 11 -> Arial 10 (legacy default font).
# ---------------------------------------------------------------------
$normalStyle = $wb.Styles.Item("Normal")
$normalStyle.Font.Name = "Arial"
$normalStyle.Font.Size = 10

Write-Host "edit applied"
